$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card11")

# --- Row 19: fill previously-blank cells with the "nan" placeholder text ---
$ws.Cells.Item(19, 2).Value  = "nan"   # B19
$ws.Cells.Item(19, 3).Value  = "nan"   # C19
$ws.Cells.Item(19, 4).Value  = "nan"   # D19
$ws.Cells.Item(19, 5).Value  = "nan"   # E19
$ws.Cells.Item(19, 6).Value  = "nan"   # F19
$ws.Cells.Item(19, 7).Value  = "nan"   # G19
$ws.Cells.Item(19, 8).Value  = "nan"   # H19
$ws.Cells.Item(19, 9).Value  = "nan"   # I19
$ws.Cells.Item(19, 10).Value = "nan"   # J19
$ws.Cells.Item(19, 11).Value = "nan"   # K19
$ws.Cells.Item(19, 13).Value = "nan"   # M19
$ws.Cells.Item(19, 16).Value = "nan"   # P19

# --- Row 20: brand-new service event for Card11 ---
$ws.Cells.Item(20, 1).Value = "'11"
$ws.Cells.Item(20, 1).Style = "Normal"
$ws.Cells.Item(20, 12).Value = "17\4\2025"
$ws.Cells.Item(20, 14).Value = "تم تغيير الجرائد الاماميه (1_2_4_5_7_8)"
$ws.Cells.Item(20, 15).Value = "الخبير"
